# Update scripts with new TPM values (Gnai2-Lhcgr, YoungD0, LR-pairs_lrc2p)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> MuSCs (Gnai2/Lhcgr) ---
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 63.91118233333333
$ws.Range("H2").Value = 191.733547
$ws.Range("I2").Value = 0.4067926910433548
$ws.Range("J2").Value = 0.4067926910433549
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.141084
$ws.Range("N2").Value = 6.423252
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 136.8392099149826
$ws.Range("R2").Value = 1231.552889234844
$ws.Range("S2").Value = 0.4067926910433548
$ws.Range("T2").Value = 0.4067926910433549

# --- Row 3: was "ECs -> MuSCs", becomes "FAPs -> MuSCs" ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 57.4434
$ws.Range("H3").Value = 172.3302
$ws.Range("I3").Value = 0.3656254573230189
$ws.Range("J3").Value = 0.365625457323019
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.141084
$ws.Range("N3").Value = 6.423252
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 122.9911446456
$ws.Range("R3").Value = 1106.9203018104
$ws.Range("S3").Value = 0.3656254573230189
$ws.Range("T3").Value = 0.365625457323019

# --- Row 4: was "FAPs -> ECs", becomes "MuSCs -> MuSCs" ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 35.755375
$ws.Range("H4").Value = 107.266125
$ws.Range("I4").Value = 0.2275818516336261
$ws.Range("J4").Value = 0.2275818516336262
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.141084
$ws.Range("N4").Value = 6.423252
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 76.55526132649999
$ws.Range("R4").Value = 688.9973519384999
$ws.Range("S4").Value = 0.2275818516336261
$ws.Range("T4").Value = 0.2275818516336262

# --- Remove the now-obsolete rows 5-7 (MuSCs target-cluster permutations) ---
$ws.Range("A5:T7").EntireRow.Delete()
